# The "IDs" list in column A was refreshed with an updated export that no
# longer drops rows early (see commit message: "fixing not going the full
# number of requests?"). The sheet grows from 65 to 75 IDs, several old
# IDs are removed, ~20 new IDs are appended, and the remaining IDs are
# reordered to match the new export order. Simplest reliable fix: replace
# the whole A1:A75 column with the corrected, final list of IDs in order.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$values = @(
    702850,703400,701350,703495,703497,703498,699640,701747,703740,703873,
    703884,703887,703888,703205,703785,703788,703792,703794,703869,703922,
    703929,703931,703935,703937,703347,702951,704031,704036,704049,704028,
    704029,703759,704061,704063,704067,704069,704077,704120,704121,704127,
    704133,704143,699827,704219,704222,704226,704239,704241,704248,704235,
    704237,703803,704259,704116,704055,704324,704418,704607,704608,704613,
    704630,704666,704672,704681,704696,704803,704805,704870,704928,704929,
    704941,704947,704949,704960,704976
)

for ($i = 0; $i -lt $values.Count; $i++) {
    $ws.Cells.Item($i + 1, 1).Value = $values[$i]
}

# Match the author's final selection/scroll state from the saved file
# (sheetView selection activeCell="B66" sqref="B66", topLeftCell="A49").
$ws.Range("B66").Select()
$excel.ActiveWindow.ScrollRow = 49
